$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 16249.75
$ws.Range("I62").Value = 27499.5
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 27499.5
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -26875.5
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 16249.75
$ws.Range("I65").Value = 27499.5
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 137497.5
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -134377.5
$ws.Range("N65").Value = -31240
$ws.Range("H86").Value = 8203.294
$ws.Range("I86").Value = 8811.875
$ws.Range("K86").Value = 8811.875
$ws.Range("M86").Value = -7688.875
$ws.Range("H89").Value = 8203.294
$ws.Range("I89").Value = 8811.875
$ws.Range("K89").Value = 44059.375
$ws.Range("M89").Value = -38443.375
$ws.Range("H97").Value = 7900
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").Value = ""
$ws.Range("H98").Value = 2406.6428
$ws.Range("I98").Value = 2406.6428
$ws.Range("K98").Value = 2406.6428
$ws.Range("M98").Value = -908.6428000000001
$ws.Range("H106").Value = 7375
$ws.Range("I106").Value = 7375
$ws.Range("K106").Value = 7375
$ws.Range("M106").Value = -6744
$ws.Range("H107").Value = 1131.4117
$ws.Range("I107").Value = 264.2143
$ws.Range("K107").Value = 264.2143
$ws.Range("M107").Value = 1655.7857
$ws.Range("H112").Value = 2440.125
$ws.Range("J112").Value = 2753.3076
$ws.Range("L112").Value = 8259.9228
$ws.Range("N112").Value = -10475.9228
$ws.Range("H122").Value = 2406.6428
$ws.Range("I122").Value = 2406.6428
$ws.Range("K122").Value = 7219.928400000001
$ws.Range("M122").Value = -4769.928400000001
$ws.Range("H137").Value = 13893255
$ws.Range("I137").Value = 23811304
$ws.Range("K137").Value = 71433912
$ws.Range("M137").Value = -71431362
$ws.Range("H138").Value = 3899.3215
$ws.Range("I138").Value = 3971.8
$ws.Range("J138").Value = 3892.2156
$ws.Range("K138").Value = 11915.4
$ws.Range("L138").Value = 11676.6468
$ws.Range("M138").Value = -6775.400000000001
$ws.Range("N138").Value = -21956.6468
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1669.7
$ws.Range("I45").Value = 1449.5
$ws.Range("K45").Value = 1449.5
$ws.Range("M45").Value = -1072.5
$ws.Range("H61").Value = 1758036.9
$ws.Range("I61").Value = 3579.5957
$ws.Range("K61").Value = 3579.5957
$ws.Range("M61").Value = -3367.5957
$ws.Range("H97").Value = 23811208
$ws.Range("I97").Value = 1100.3513
$ws.Range("J97").Value = 200006000
$ws.Range("K97").Value = 1100.3513
$ws.Range("L97").Value = 200006000
$ws.Range("M97").Value = -604.3513
$ws.Range("N97").Value = -200006992
$ws.Range("H132").Value = 628201.2
$ws.Range("I132").Value = 738432.4
$ws.Range("K132").Value = 2215297.2
$ws.Range("M132").Value = -2212767.2
$ws.Range("H136").Value = 1758036.9
$ws.Range("I136").Value = 3579.5957
$ws.Range("K136").Value = 10738.7871
$ws.Range("M136").Value = -8188.7871
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 102499.5
$ws.Range("J58").Value = 102499.5
$ws.Range("L58").Value = 102499.5
$ws.Range("N58").Value = -103087.5
$ws.Range("H134").Value = 8785692
$ws.Range("I134").Value = 10407.728
$ws.Range("K134").Value = 31223.184
$ws.Range("M134").Value = -28688.184
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2528744.5
$ws.Range("I31").Value = 3705703.5
$ws.Range("K31").Value = 3705703.5
$ws.Range("M31").Value = -3705408.5
$ws.Range("H34").Value = 2528744.5
$ws.Range("I34").Value = 3705703.5
$ws.Range("K34").Value = 3705703.5
$ws.Range("M34").Value = -3705501.5
$ws.Range("H58").Value = 2882097.5
$ws.Range("I58").Value = 4621
$ws.Range("K58").Value = 4621
$ws.Range("M58").Value = -4418
$ws.Range("H99").Value = 32803.5
$ws.Range("I99").Value = 33174
$ws.Range("J99").Value = 31903.715
$ws.Range("K99").Value = 33174
$ws.Range("L99").Value = 31903.715
$ws.Range("M99").Value = -31676
$ws.Range("N99").Value = -34899.715
$ws.Range("H122").Value = 6669.121
$ws.Range("J122").Value = 82834.5
$ws.Range("L122").Value = 248503.5
$ws.Range("N122").Value = -253403.5
$ws.Range("H126").Value = 32803.5
$ws.Range("I126").Value = 33174
$ws.Range("J126").Value = 31903.715
$ws.Range("K126").Value = 99522
$ws.Range("L126").Value = 95711.145
$ws.Range("M126").Value = -97052
$ws.Range("N126").Value = -100651.145
$ws.Range("H136").Value = 2882097.5
$ws.Range("I136").Value = 4621
$ws.Range("K136").Value = 13863
$ws.Range("M136").Value = -11313
$ws.Range("H141").Value = 234462.8
$ws.Range("J141").Value = 255149.47
$ws.Range("L141").Value = 255149.47
$ws.Range("N141").Value = -265509.47
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 399.5
$ws.Range("I8").Value = 399.5
$ws.Range("K8").Value = 1198.5
$ws.Range("M8").Value = -1059.5
$ws.Range("H57").Value = 900
$ws.Range("I57").Value = 900
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 2700
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -2141
$ws.Range("N57").Value = ""
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 11634.818
$ws.Range("H46").Value = 3283
$ws.Range("I46").Value = 3283
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 3283
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -3127
$ws.Range("N46").Value = ""
$ws.Range("H57").Value = 14439.857
$ws.Range("J57").Value = 29998.5
$ws.Range("L57").Value = 29998.5
$ws.Range("N57").Value = -31638.5
$ws.Range("H70").Value = 22751.727
$ws.Range("I70").Value = 25458.6
$ws.Range("J70").Value = 20496
$ws.Range("K70").Value = 25458.6
$ws.Range("L70").Value = 20496
$ws.Range("M70").Value = -25188.6
$ws.Range("N70").Value = -21036
$ws.Range("H73").Value = 22751.727
$ws.Range("I73").Value = 25458.6
$ws.Range("J73").Value = 20496
$ws.Range("K73").Value = 25458.6
$ws.Range("L73").Value = 20496
$ws.Range("M73").Value = -24522.6
$ws.Range("N73").Value = -22368
$ws.Range("H80").Value = 7473.75
$ws.Range("J80").Value = 4998.5
$ws.Range("L80").Value = 4998.5
$ws.Range("N80").Value = -6994.5
$ws.Range("H83").Value = 7473.75
$ws.Range("J83").Value = 4998.5
$ws.Range("M83").Value = -34979
$ws.Range("N83").Value = -34976.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6963.7144
$ws.Range("I40").Value = 6730.154
$ws.Range("K40").Value = 6730.154
$ws.Range("M40").Value = -6594.154
$ws.Range("H74").Value = 41233.25
$ws.Range("I74").Value = 41233.25
$ws.Range("K74").Value = 41233.25
$ws.Range("M74").Value = -40235.25
$ws.Range("H77").Value = 41233.25
$ws.Range("I77").Value = 41233.25
$ws.Range("K77").Value = 123699.75
$ws.Range("M77").Value = -118707.75
$ws.Range("H140").Value = 53666
$ws.Range("J140").Value = 57999
$ws.Range("L140").Value = 57999
$ws.Range("N140").Value = -68359
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 91462.53999999999
$ws.Range("I122").Value = 6716.5557
$ws.Range("J122").Value = 282141
$ws.Range("K122").Value = 20149.6671
$ws.Range("L122").Value = 846423
$ws.Range("M122").Value = -17699.6671
$ws.Range("N122").Value = -851323
$ws.Range("H126").Value = 3299.6
$ws.Range("I126").Value = 3374.5
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 10123.5
$ws.Range("L126").Value = 9000
$ws.Range("N126").Value = -13940
$ws.Range("M126").Value = -7653.5
$ws.Range("H140").Value = 107857.4
$ws.Range("J140").Value = 107857.4
$ws.Range("L140").Value = 107857.4
$ws.Range("N140").Value = -118217.4
